$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 163413
$ws.Range("C4").Value = 154404
$ws.Range("C5").Value = 9009
$ws.Range("C7").Value = 5.51
$ws.Range("C8").Value = 64.73999999999999
